# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the newly generated site data (gh-pages build 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9228
    $ws.Range("F4").Value = 488
}
